# Daily attendance processing - 2026-01-27 20:38:15
#
# Normalize the "Recorded By" values (column G) on the session-analysis
# sheet: rows that were recorded by both "System" and the user are
# currently listed as "System, dnasr281@gmail.com"; flip the ordering so
# the human user is listed first: "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$searchRange = $ws.UsedRange

$firstFound = $searchRange.Find($target)
if ($firstFound -ne $null) {
    $firstAddress = $firstFound.Address()
    $cell = $firstFound
    $keepGoing = $true
    while ($keepGoing) {
        $cell.Value = $replacement
        $cell = $searchRange.FindNext($cell)
        if ($cell -eq $null -or $cell.Address() -eq $firstAddress) {
            $keepGoing = $false
        }
    }
}
